$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert first new weekly record at row 346 ---
$ws.Rows.Item(346).Insert()

$ws.Range("A346").Value = 5
$ws.Range("B346").Value = "Macroferia Regional de Talca"
$ws.Range("C346").Value = "Maule"
$ws.Range("D346").Value = 45120
$ws.Range("E346").Value = 7
$ws.Range("F346").Value = 100112009
$ws.Range("G346").Value = "Acelga"
$ws.Range("H346").Value = "Sin especificar"
$ws.Range("I346").Value = "Primera"
$ws.Range("J346").Value = 500
$ws.Range("K346").Value = 1800
$ws.Range("L346").Value = 1800
$ws.Range("M346").Value = 1800
$ws.Range("N346").Value = "$/docena de atados (4 kilos)"
$ws.Range("O346").Value = "Región del Maule"
$ws.Range("P346").Value = 450
$ws.Range("Q346").Value = 4
$ws.Range("R346").Value = "Hortaliza"

# --- Insert second new weekly record at row 462 ---
$ws.Rows.Item(462).Insert()

$ws.Range("A462").Value = 5
$ws.Range("B462").Value = "Macroferia Regional de Talca"
$ws.Range("C462").Value = "Maule"
$ws.Range("D462").Value = 45121
$ws.Range("E462").Value = 7
$ws.Range("F462").Value = 100112009
$ws.Range("G462").Value = "Acelga"
$ws.Range("H462").Value = "Sin especificar"
$ws.Range("I462").Value = "Primera"
$ws.Range("J462").Value = 500
$ws.Range("K462").Value = 1800
$ws.Range("L462").Value = 1800
$ws.Range("M462").Value = 1800
$ws.Range("N462").Value = "$/docena de atados (4 kilos)"
$ws.Range("O462").Value = "Región del Maule"
$ws.Range("P462").Value = 450
$ws.Range("Q462").Value = 4
$ws.Range("R462").Value = "Hortaliza"
